$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 31.27132566666667
$ws.Range("H2").Value = 93.81397700000001
$ws.Range("I2").Value = 0.9493361071405608
$ws.Range("J2").Value = 0.9493361071405608
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.397026
$ws.Range("N2").Value = 64.191078
$ws.Range("O2").Value = 0.7694665596935515
$ws.Range("P2").Value = 0.7694665596935515
$ws.Range("Q2").Value = 669.1133683441341
$ws.Range("R2").Value = 6022.020315097207
$ws.Range("S2").Value = 0.7304823883543161
$ws.Range("T2").Value = 0.7304823883543161

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 31.27132566666667
$ws.Range("H3").Value = 93.81397700000001
$ws.Range("I3").Value = 0.9493361071405608
$ws.Range("J3").Value = 0.9493361071405608
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.818642333333333
$ws.Range("N3").Value = 14.455927
$ws.Range("O3").Value = 0.1732850228168956
$ws.Range("P3").Value = 0.1732850228168955
$ws.Range("Q3").Value = 150.6853336768532
$ws.Range("R3").Value = 1356.168003091679
$ws.Range("S3").Value = 0.1645057289867549
$ws.Range("T3").Value = 0.1645057289867549

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 31.27132566666667
$ws.Range("H4").Value = 93.81397700000001
$ws.Range("I4").Value = 0.9493361071405608
$ws.Range("J4").Value = 0.9493361071405608
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.591941666666667
$ws.Range("N4").Value = 4.775825
$ws.Range("O4").Value = 0.05724841748955292
$ws.Range("P4").Value = 0.05724841748955292
$ws.Range("Q4").Value = 49.78212630066945
$ws.Range("R4").Value = 448.039136706025
$ws.Range("S4").Value = 0.05434798979948977
$ws.Range("T4").Value = 0.05434798979948977

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.96805
$ws.Range("H5").Value = 2.90415
$ws.Range("I5").Value = 0.02938809912676722
$ws.Range("J5").Value = 0.02938809912676721
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.397026
$ws.Range("N5").Value = 64.191078
$ws.Range("O5").Value = 0.7694665596935515
$ws.Range("P5").Value = 0.7694665596935515
$ws.Range("Q5").Value = 20.7133910193
$ws.Range("R5").Value = 186.4205191737
$ws.Range("S5").Value = 0.02261315953100664
$ws.Range("T5").Value = 0.02261315953100663

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.96805
$ws.Range("H6").Value = 2.90415
$ws.Range("I6").Value = 0.02938809912676722
$ws.Range("J6").Value = 0.02938809912676721
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.818642333333333
$ws.Range("N6").Value = 14.455927
$ws.Range("O6").Value = 0.1732850228168956
$ws.Range("P6").Value = 0.1732850228168955
$ws.Range("Q6").Value = 4.664686710783333
$ws.Range("R6").Value = 41.98218039704999
$ws.Range("S6").Value = 0.005092517427727045
$ws.Range("T6").Value = 0.005092517427727044

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.96805
$ws.Range("H7").Value = 2.90415
$ws.Range("I7").Value = 0.02938809912676722
$ws.Range("J7").Value = 0.02938809912676721
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.591941666666667
$ws.Range("N7").Value = 4.775825
$ws.Range("O7").Value = 0.05724841748955292
$ws.Range("P7").Value = 0.05724841748955292
$ws.Range("Q7").Value = 1.541079130416667
$ws.Range("R7").Value = 13.86971217375
$ws.Range("S7").Value = 0.001682422168033535
$ws.Range("T7").Value = 0.001682422168033535

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.700829
$ws.Range("H8").Value = 2.102487
$ws.Range("I8").Value = 0.02127579373267201
$ws.Range("J8").Value = 0.02127579373267201
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.397026
$ws.Range("N8").Value = 64.191078
$ws.Range("O8").Value = 0.7694665596935515
$ws.Range("P8").Value = 0.7694665596935515
$ws.Range("Q8").Value = 14.995656334554
$ws.Range("R8").Value = 134.960907010986
$ws.Range("S8").Value = 0.01637101180822876
$ws.Range("T8").Value = 0.01637101180822876

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.700829
$ws.Range("H9").Value = 2.102487
$ws.Range("I9").Value = 0.02127579373267201
$ws.Range("J9").Value = 0.02127579373267201
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.818642333333333
$ws.Range("N9").Value = 14.455927
$ws.Range("O9").Value = 0.1732850228168956
$ws.Range("P9").Value = 0.1732850228168955
$ws.Range("Q9").Value = 3.377044287827667
$ws.Range("R9").Value = 30.393398590449
$ws.Range("S9").Value = 0.003686776402413633
$ws.Range("T9").Value = 0.003686776402413632

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.700829
$ws.Range("H10").Value = 2.102487
$ws.Range("I10").Value = 0.02127579373267201
$ws.Range("J10").Value = 0.02127579373267201
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.591941666666667
$ws.Range("N10").Value = 4.775825
$ws.Range("O10").Value = 0.05724841748955292
$ws.Range("P10").Value = 0.05724841748955292
$ws.Range("Q10").Value = 1.115678886308334
$ws.Range("R10").Value = 10.041109976775
$ws.Range("S10").Value = 0.001218005522029621
$ws.Range("T10").Value = 0.001218005522029621

